$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the summary-statistics rows (Count/Mean/Stdev/Min/Max labels + values)
# that were computed over the raw data block; keep the cell styling intact.
$ws.Range("E51:AO60").ClearContents()

# Update the active selection / view to match the author's final state.
[void]$ws.Range("D49:AQ68").Select()
